$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Two new columns are appended to the player table:
#   AF = "Offcourt"     -> the player's TeamTS value from the "paired"
#                          On/Off row (mirrors column AD, 19 rows down
#                          for rows 3-21, 19 rows up for rows 22-40)
#   AG = "Differential"  -> AD(row) - AD(paired row)
# ---------------------------------------------------------------------

# --- Headers -----------------------------------------------------------
# Copy the header cell formatting (centered, merged look) from AE1:AE2,
# which is the column immediately to the left of the new ones.
$ws.Range("AE1:AE2").Copy() | Out-Null
$ws.Range("AF1:AF2").PasteSpecial(-4122) | Out-Null
$ws.Range("AG1:AG2").PasteSpecial(-4122) | Out-Null

$ws.Range("AF1").Value = "Offcourt"
$ws.Range("AG1").Value = "Differential"

$ws.Range("AF1:AF2").Merge() | Out-Null
$ws.Range("AG1:AG2").Merge() | Out-Null

# --- Value-column formatting --------------------------------------------
# Column AF should look like column AD (right aligned, alternating banded
# shading, same borders) so copy that formatting across the data rows.
$ws.Range("AD3:AD40").Copy() | Out-Null
$ws.Range("AF3:AF40").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- AF: mirror of the paired row's TeamTS (AD) value -------------------
for ($r = 3; $r -le 21; $r++) {
    $pair = $r + 19
    $ws.Range("AF$r").Value = $ws.Range("AD$pair").Value2
}
for ($r = 22; $r -le 40; $r++) {
    $pair = $r - 19
    $ws.Range("AF$r").Value = $ws.Range("AD$pair").Value2
}

# --- AG: differential between the row's own AD and its paired AD --------
# Top half (rows 3-21) keeps live formulas; row 3 / row 21 are the two
# boundary rows (standalone formulas), rows 4-20 share the same pattern.
$ws.Range("AG3").Formula = "=AD3-AD22"
for ($r = 4; $r -le 20; $r++) {
    $pair = $r + 19
    $ws.Range("AG$r").Formula = "=AD$r-AD$pair"
}
$ws.Range("AG21").Formula = "=AD21-AD40"

# Bottom half (rows 22-40) stores the same differential as a plain value
# (no formula), mirroring the top half's numbers.
for ($r = 22; $r -le 40; $r++) {
    $pair = $r - 19
    $ws.Range("AG$r").Value = $ws.Range("AD$pair").Value2 - $ws.Range("AD$r").Value2
}

# --- View bookkeeping -----------------------------------------------------
# Scroll the sheet so column B is the left-most visible column and move
# the active selection onto the new trailing column.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AH21").Select() | Out-Null

Write-Host "done"
